$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of resale-number data for 2023-05-30 14:05:52 (row 8)
$rowIndex = 8

$textValues = @{
    "A" = "2023-05-30"
    "B" = "14:05:52"
    "C" = "Tuesday"
    "D" = "22"
}

$numberValues = @{
    "E" = 119835
    "F" = 133541
    "G" = 158282
    "H" = 130822
    "I" = 174497
    "J" = 113620
    "K" = 198303
    "L" = 220206
    "M" = 171994
    "N" = 119875
    "O" = 38604
    "P" = 34861
    "Q" = 50503
    "R" = -1
    "S" = 36734
    "T" = -1
}

# Mark the text columns as Text format first so that date/time/number-looking
# strings ("2023-05-30", "14:05:52", "22") are stored as plain text instead of
# being auto-converted into date/time/number values.
foreach ($col in $textValues.Keys) {
    $ws.Range("$col$rowIndex").NumberFormat = "@"
}

foreach ($col in $textValues.Keys) {
    $ws.Range("$col$rowIndex").Value2 = $textValues[$col]
}

# Remove the temporary text formatting again so the new cells end up without
# an explicit style, matching the rest of the data rows.
foreach ($col in $textValues.Keys) {
    $ws.Range("$col$rowIndex").ClearFormats()
}

foreach ($col in $numberValues.Keys) {
    $ws.Range("$col$rowIndex").Value2 = $numberValues[$col]
}
